# act tablas web jul25
# Renames the sheets (Datos -> Data, Ficha técnica -> Metadata), refreshes
# the "Data" time series (now 2007-2024, newest first, 2020/2021 omitted)
# and rewrites the "Metadata" fact sheet with lower-cased field keys and an
# updated citation / new "observaciones" + "actualizacion" rows.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# --- rename sheets -----------------------------------------------------
$ws1.Name = "Data"
$ws2.Name = "Metadata"

# --- wipe existing content so the shared-string table gets rebuilt -----
# clean (only strings actually referenced survive) and in the exact order
# we (re)write the cells below.
$ws1.Cells.Clear()
$ws2.Cells.Clear()

# --- helper: write a cell as literal text (not auto-coerced to a number)
function Set-TextCell($cell, $text) {
    $cell.Value = "'" + $text
    $cell.Style = "Normal"
}

# --- Data sheet ----------------------------------------------------------
$ws1.Cells.Item(1,1).Value = "Fecha"
$ws1.Cells.Item(1,2).Value = "Valor"

$data = @(
    @("2024", 5.5),
    @("2023", 5.8),
    @("2022", 4.6),
    @("2019", 4.5),
    @("2018", 4.7),
    @("2017", 4.7),
    @("2016", 5),
    @("2015", 5.1),
    @("2014", 5),
    @("2013", 4.8),
    @("2012", 5.5),
    @("2011", 5),
    @("2010", 5.9),
    @("2009", 5.4),
    @("2008", 5.6),
    @("2007", 5.5)
)

$r = 2
foreach ($row in $data) {
    Set-TextCell $ws1.Cells.Item($r,1) $row[0]
    $ws1.Cells.Item($r,2).Value = $row[1]
    $r = $r + 1
}

# --- Metadata sheet ------------------------------------------------------
$observaciones = "Desde marzo de 2020 hasta junio de 2021 se interrumpió el relevamiento presencial y se aplicó de manera telefónica un cuestionario restringido con el objetivo de continuar publicando los indicadores de ingresos y mercado de trabajo. En ese período la encuesta pasó a ser de paneles rotativos elegidos al azar a partir de los casos respondentes del año anterior. `nEn julio de 2021 el INE retomó la realización de encuestas presenciales, pero introdujo un cambio metodológico, ya que la ECH pasa a ser una encuesta de panel rotativo con periodicidad mensual compuesta por seis paneles o grupos de rotación, cada uno de los cuales es una muestra representativa de la población. Con esta nueva metodología, cada hogar seleccionado participa durante seis meses de la ECH."

$cita = "UMAD con base en Instituto de Economía, Universidad de la República (2020) Encuesta Continua de Hogares Compatibilizada 1981-2018 Versión 12 DOI: http://doiorg/1047426/ECHINE (Hasta 2019) / A partir de 2020 con base en ECH - INE`n"

$meta = @(
    @(" ", " "),
    @("nomindicador", "Porcentaje de personas que viven en asentamientos"),
    @("derecho", "Vivienda"),
    @("conindicador", "Asentamientos"),
    @("tipoind", "Resultados"),
    @("definicion", "El indicador mide el porcentaje de personas en viviendas ubicadas en asentamiento irregular."),
    @("calculo", "Para cada año calcular: (Cantidad de personas que residen en viviendas ubicadas en asentamiento irregular / Cantidad total de personas en viviendas particulares)*100"),
    @("observaciones", $observaciones),
    @("actualizacion", "Julio 2025"),
    @("cita", $cita),
    @("Mirador DESCA - UMAD/FCS – INDDHH", " ")
)

$r = 1
foreach ($row in $meta) {
    $ws2.Cells.Item($r,1).Value = $row[0]
    $ws2.Cells.Item($r,2).Value = $row[1]
    $r = $r + 1
}
